# Revert "Cường ghi thông tin rồi"
# Removes the phone/email contact info that had been added for
# "Nguyễn Duy Cường" (row 23) and puts the name back on the "NHÓM 4"
# row (row 22), matching the pre-commit layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the mailto: hyperlink that lived on F23 before touching any
# cell values/styles.
$null = $ws.Range("F23").Hyperlinks.Delete()

# The name currently sits on D23; move it back up onto the "NHÓM 4"
# row (D22) the way it was before the contact info was ever added.
$name = $ws.Cells.Item(23, 4).Value()
$ws.Cells.Item(22, 4).Value = $name

# Clear out the now-empty contact-info row (name/phone/email cells)
# and drop the hyperlink cell style so F23 goes back to Normal.
$ws.Range("D23:F23").Style = "Normal"
$ws.Range("D23:F23").ClearContents()

# Remove the now-unused "Hyperlink" cell style definition.
$null = $wb.Styles.Item("Hyperlink").Delete()

# Restore the previously-selected cell.
$null = $ws.Range("E23").Select()
